$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row after the current row 12 (old A12=45805 row),
#        pushing old rows 13..35 down to 14..36. This produces the new
#        blank row 13 and keeps dates realigned (new A14 = 45806, etc.)
$ws.Rows("13:13").Insert()

# --- 2. Fill in the "1h" / activity text for the (now augmented) row 12
$ws.Range("B12").Value = "1h"
$ws.Range("C12").Value = "Gestion stage & gestion recherche de stage"

# Give row 12 the same date-cell style used by the rest of the A column
# block above it (s=11) so it matches the row that now carries a
# secondary "sub row" (A13) beneath it, and restore its 30pt height.
$ws.Range("A12").RowHeight = 30

# --- 3. Merge A12:A13 (the date cell now spans the new blank sub-row)
$ws.Range("A12:A13").Merge()

# --- 4. Row 2 no longer needs an explicit 30pt override; let it
#        auto-size back down to the sheet's default row height.
$ws.Rows("2:2").AutoFit()

# --- 5. Widen column C slightly to fit the new text.
$ws.Columns("C:C").ColumnWidth = 34.45

# --- 6. Move the active selection to B13 (matches the author's cursor
#        position after the edit).
$ws.Range("B13").Select()
